$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 359.7143
$ws.Range("I33").Value = 209.75
$ws.Range("J33").Value = 839.6
$ws.Range("K33").Value = 209.75
$ws.Range("L33").Value = 839.6
$ws.Range("M33").Value = 19.25
$ws.Range("N33").Value = -1297.6
# Row 40
$ws.Range("H40").Value = 4993.154
$ws.Range("I40").Value = 4700
$ws.Range("J40").Value = 5123.4443
$ws.Range("K40").Value = 4700
$ws.Range("L40").Value = 5123.4443
$ws.Range("M40").Value = -4525
$ws.Range("N40").Value = -5473.4443
# Row 100
$ws.Range("H100").Value = 5656.4546
$ws.Range("I100").Value = 4735.8
$ws.Range("K100").Value = 4735.8
$ws.Range("M100").Value = -4194.8
# Row 116
$ws.Range("H116").Value = 10000.25
$ws.Range("J116").Value = 9666.5
$ws.Range("L116").Value = 9666.5
$ws.Range("N116").Value = -16550.5
# Row 123
$ws.Range("H123").Value = 49166.668
$ws.Range("J123").Value = 49166.668
$ws.Range("L123").Value = 49166.668
$ws.Range("N123").Value = -58966.668
# Row 127
$ws.Range("H127").Value = 12509.272
$ws.Range("I127").Value = 1371.4286
$ws.Range("K127").Value = 4114.2858
$ws.Range("M127").Value = 845.7142000000003
# Row 132
$ws.Range("H132").Value = 16648.71
$ws.Range("I132").Value = 2700.08
$ws.Range("K132").Value = 8100.24
$ws.Range("M132").Value = -5570.24
# Row 135
$ws.Range("H135").Value = 22729520
$ws.Range("I135").Value = 26317972
$ws.Range("K135").Value = 236861748
$ws.Range("M135").Value = -236859213

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2843.6843
$ws.Range("I45").Value = 1337.9286
$ws.Range("J45").Value = 7059.8
$ws.Range("K45").Value = 1337.9286
$ws.Range("L45").Value = 7059.8
$ws.Range("M45").Value = -960.9286
$ws.Range("N45").Value = -7813.8
# Row 110
$ws.Range("H110").Value = 3779
$ws.Range("I110").Value = 3612.85
$ws.Range("K110").Value = 3612.85
$ws.Range("M110").Value = -1567.85
# Row 122
$ws.Range("H122").Value = 4744.4517
$ws.Range("I122").Value = 4740.0835
$ws.Range("K122").Value = 14220.2505
$ws.Range("M122").Value = -11770.2505
# Row 139
$ws.Range("H139").Value = 79818.17999999999
$ws.Range("J139").Value = 79818.17999999999
$ws.Range("L139").Value = 79818.17999999999
$ws.Range("N139").Value = -90098.17999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3969.75
$ws.Range("I20").Value = 2695.75
$ws.Range("K20").Value = 2695.75
$ws.Range("M20").Value = -2448.75
# Row 99
$ws.Range("H99").Value = 2676.7778
$ws.Range("I99").Value = 2811.625
$ws.Range("K99").Value = 2811.625
$ws.Range("M99").Value = -1313.625
# Row 105
$ws.Range("H105").Value = 3074.1667
$ws.Range("I105").Value = 1618.8182
$ws.Range("J105").Value = 5361.143
$ws.Range("K105").Value = 1618.8182
$ws.Range("L105").Value = 5361.143
$ws.Range("M105").Value = 128.1818000000001
$ws.Range("N105").Value = -8855.143
# Row 107
$ws.Range("H107").Value = 4658.857
$ws.Range("J107").Value = 6443.4
$ws.Range("L107").Value = 6443.4
$ws.Range("N107").Value = -10283.4
# Row 134
$ws.Range("H134").Value = 3060.3684
$ws.Range("I134").Value = 2540.2307
$ws.Range("J134").Value = 4187.3335
$ws.Range("K134").Value = 7620.6921
$ws.Range("L134").Value = 12562.0005
$ws.Range("M134").Value = -5085.6921
$ws.Range("N134").Value = -17632.0005

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2460.25
$ws.Range("I16").Value = 2227.45
$ws.Range("K16").Value = 2227.45
$ws.Range("M16").Value = -1940.45
# Row 25
$ws.Range("H25").Value = 543
$ws.Range("I25").Value = 543
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 543
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -369
$ws.Range("N25").ClearContents()
# Row 58
$ws.Range("H58").Value = 2009.1666
$ws.Range("I58").Value = 1283.5714
$ws.Range("K58").Value = 1283.5714
$ws.Range("M58").Value = -1080.5714
# Row 113
$ws.Range("H113").Value = 2460.25
$ws.Range("I113").Value = 2227.45
$ws.Range("K113").Value = 2227.45
$ws.Range("M113").Value = -57.44999999999982
# Row 122
$ws.Range("H122").Value = 4207.3335
$ws.Range("I122").Value = 3554.8
$ws.Range("K122").Value = 10664.4
$ws.Range("M122").Value = -8214.400000000001
# Row 136
$ws.Range("H136").Value = 2009.1666
$ws.Range("I136").Value = 1283.5714
$ws.Range("K136").Value = 3850.7142
$ws.Range("M136").Value = -1300.7142

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 56
$ws.Range("H56").Value = 7230.125
$ws.Range("I56").Value = 7230.125
$ws.Range("K56").Value = 7230.125
$ws.Range("M56").Value = -6700.125
# Row 107
$ws.Range("H107").Value = 2588.4783
$ws.Range("J107").Value = 2709.6428
$ws.Range("L107").Value = 8128.928400000001
$ws.Range("N107").Value = -11968.9284

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 18
$ws.Range("H18").Value = 4900
$ws.Range("I18").Value = 4900
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 4900
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -4607
$ws.Range("N18").ClearContents()
# Row 70
$ws.Range("H70").Value = 143838.75
$ws.Range("I70").Value = 226433
$ws.Range("J70").Value = 6181.6665
$ws.Range("K70").Value = 226433
$ws.Range("L70").Value = 6181.6665
$ws.Range("M70").Value = -226163
$ws.Range("N70").Value = -6721.6665
# Row 73
$ws.Range("H73").Value = 143838.75
$ws.Range("I73").Value = 226433
$ws.Range("J73").Value = 6181.6665
$ws.Range("K73").Value = 226433
$ws.Range("L73").Value = 6181.6665
$ws.Range("M73").Value = -225497
$ws.Range("N73").Value = -8053.6665
# Row 97
$ws.Range("H97").Value = 4119.9644
$ws.Range("I97").Value = 501.1579
$ws.Range("J97").Value = 11759.667
$ws.Range("K97").Value = 501.1579
$ws.Range("L97").Value = 11759.667
$ws.Range("M97").Value = -5.157899999999984
$ws.Range("N97").Value = -12751.667
# Row 122
$ws.Range("H122").Value = 3279.7222
$ws.Range("I122").Value = 2724
$ws.Range("K122").Value = 8172
$ws.Range("M122").Value = -5722
# Row 126
$ws.Range("H126").Value = 4118
$ws.Range("I126").Value = 4257.75
$ws.Range("K126").Value = 12773.25
$ws.Range("M126").Value = -10303.25
# Row 132
$ws.Range("H132").Value = 6841.5654
$ws.Range("I132").Value = 6660.2666
$ws.Range("K132").Value = 19980.7998
$ws.Range("M132").Value = -17450.7998

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 8333.333000000001
$ws.Range("I7").Value = 11000
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 11000
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -10888
$ws.Range("N7").Value = -3224
# Row 23
$ws.Range("H23").Value = 29995
$ws.Range("I23").Value = 29995
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 29995
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -29765
$ws.Range("N23").ClearContents()
# Row 46
$ws.Range("H46").Value = 2751.5454
$ws.Range("I46").Value = 1400
$ws.Range("J46").Value = 3523.8572
$ws.Range("K46").Value = 1400
$ws.Range("L46").Value = 3523.8572
$ws.Range("M46").Value = -1212
$ws.Range("N46").Value = -3899.8572
# Row 51
$ws.Range("H51").Value = 39770.75
$ws.Range("J51").Value = 39770.75
$ws.Range("L51").Value = 39770.75
$ws.Range("N51").Value = -40726.75
# Row 61
$ws.Range("H61").Value = 2589.6667
$ws.Range("I61").Value = 2552.4546
$ws.Range("K61").Value = 2552.4546
$ws.Range("M61").Value = -2350.4546
# Row 68
$ws.Range("H68").Value = 3089.3
$ws.Range("I68").Value = 2714
$ws.Range("J68").Value = 3965
$ws.Range("K68").Value = 2714
$ws.Range("L68").Value = 3965
$ws.Range("M68").Value = -1965
$ws.Range("N68").Value = -5463
# Row 71
$ws.Range("H71").Value = 3089.3
$ws.Range("I71").Value = 2714
$ws.Range("J71").Value = 3965
$ws.Range("K71").Value = 13570
$ws.Range("L71").Value = 19825
$ws.Range("M71").Value = -9826
$ws.Range("N71").Value = -27313
# Row 100
$ws.Range("H100").Value = 91084.30499999999
$ws.Range("I100").Value = 225019.2
$ws.Range("K100").Value = 225019.2
$ws.Range("M100").Value = -224478.2
# Row 113
$ws.Range("H113").Value = 2589.6667
$ws.Range("I113").Value = 2552.4546
$ws.Range("K113").Value = 2552.4546
$ws.Range("M113").Value = -382.4546
# Row 126
$ws.Range("H126").Value = 8333.333000000001
$ws.Range("I126").Value = 11000
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 33000
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -30530
$ws.Range("N126").Value = -13940
# Row 132
$ws.Range("H132").Value = 1942.9062
$ws.Range("I132").Value = 1683
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 5049
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -2519
$ws.Range("N132").Value = -35060
# Row 136
$ws.Range("H136").Value = 3658
$ws.Range("I136").Value = 3424.6875
$ws.Range("K136").Value = 10274.0625
$ws.Range("M136").Value = -7724.0625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4866.3335
$ws.Range("J62").Value = 4866.3335
$ws.Range("L62").Value = 4866.3335
$ws.Range("N62").Value = -6114.3335
# Row 65
$ws.Range("H65").Value = 4866.3335
$ws.Range("J65").Value = 4866.3335
$ws.Range("L65").Value = 24331.6675
$ws.Range("N65").Value = -30571.6675
# Row 100
$ws.Range("H100").Value = 1169.9
$ws.Range("I100").Value = 1458.3334
$ws.Range("K100").Value = 2916.6668
$ws.Range("M100").Value = -2375.6668
# Row 136
$ws.Range("H136").Value = 1564.5714
$ws.Range("I136").Value = 1292
$ws.Range("J136").Value = 3200
$ws.Range("K136").Value = 3876
$ws.Range("L136").Value = 9600
$ws.Range("M136").Value = -1326
$ws.Range("N136").Value = -14700
